$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B2").Value = 0.50187265917603
$ws.Range("C2").Value = 0.50093808630394
$ws.Range("E2").Value = 0.6675
$ws.Range("F2").Value = 0.8338538413491567
$ws.Range("G2").Value = 0.9630965593784684
$ws.Range("H2").Value = 0.709769389386862
$ws.Range("J2").Value = 532
$ws.Range("K2").Value = 2

# --- Classification Report sheet ---
$ws = $wb.Worksheets.Item("Classification Report")
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0.003745318352059925
$ws.Range("D2").Value = 0.007462686567164179

$ws.Range("B3").Value = 0.50093808630394
$ws.Range("D3").Value = 0.6675

$ws.Range("B4").Value = 0.50187265917603
$ws.Range("C4").Value = 0.50187265917603
$ws.Range("D4").Value = 0.50187265917603
$ws.Range("E4").Value = 0.50187265917603

$ws.Range("B5").Value = 0.75046904315197
$ws.Range("C5").Value = 0.50187265917603
$ws.Range("D5").Value = 0.3374813432835821

$ws.Range("B6").Value = 0.7504690431519699
$ws.Range("C6").Value = 0.50187265917603
$ws.Range("D6").Value = 0.3374813432835821

# --- Confusion Matrix sheet ---
$ws = $wb.Worksheets.Item("Confusion Matrix")
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 532
